# The footer's first paragraph currently holds a centered "-<PAGE>-" page
# number field (built from a literal "-" run, an SDT-wrapped PAGE field,
# and a trailing literal "-" run). The edit removes that whole page-number
# construct, leaving a single empty "Footer"-styled paragraph (no center
# justification either) - i.e. the footer becomes blank.

$d = $word.ActiveDocument

$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)  # wdHeaderFooterPrimary

if ($ftr.Range.Paragraphs.Count -gt 1) {
    # Deleting the first paragraph's Range (including its end-of-paragraph
    # mark) removes the "-", the PAGE-number SDT/field, and the trailing
    # "-", and merges what follows up to the first position - which is the
    # already-empty, non-centered "Footer" paragraph that trails it. That
    # leaves exactly one empty paragraph carrying just the Footer style.
    $p1 = $ftr.Range.Paragraphs.Item(1)
    $p1.Range.Delete()
}
elseif ($ftr.Range.Text -ne "") {
    # Fallback in case only a single footer paragraph exists: clear its
    # text/fields and drop the centered alignment directly.
    $ftr.Range.Fields | ForEach-Object { $_.Delete() }
    $ftr.Range.Text = ""
    $ftr.Range.ParagraphFormat.Alignment = 0
}

Write-Output ("Footer text after edit: [" + $ftr.Range.Text + "]")
